$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1437.0938
$ws.Range("J129").Value = 1570.9642
$ws.Range("L129").Value = 4712.892599999999
$ws.Range("N129").Value = -14712.8926
$ws.Range("H138").Value = 5036.155
$ws.Range("I138").Value = 8939.4
$ws.Range("J138").Value = 4740.4546
$ws.Range("K138").Value = 26818.2
$ws.Range("L138").Value = 14221.3638
$ws.Range("M138").Value = -21678.2
$ws.Range("N138").Value = -24501.3638
$ws.Range("H140").Value = 76497.5
$ws.Range("J140").Value = 76497.5
$ws.Range("L140").Value = 76497.5
$ws.Range("N140").Value = -86857.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21757454
$ws.Range("I32").Value = 23826306
$ws.Range("K32").Value = 23826306
$ws.Range("M32").Value = -23826019
$ws.Range("H45").Value = 2003.5555
$ws.Range("I45").Value = 1037.5
$ws.Range("J45").Value = 2776.4
$ws.Range("K45").Value = 1037.5
$ws.Range("L45").Value = 2776.4
$ws.Range("M45").Value = -660.5
$ws.Range("N45").Value = -3530.4
$ws.Range("H63").Value = 4711.3657
$ws.Range("I63").Value = 3376.6667
$ws.Range("J63").Value = 5086.75
$ws.Range("K63").Value = 3376.6667
$ws.Range("L63").Value = 5086.75
$ws.Range("M63").Value = -2690.6667
$ws.Range("N63").Value = -6458.75
$ws.Range("H64").Value = 28000
$ws.Range("J64").Value = 28000
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496
$ws.Range("H66").Value = 4711.3657
$ws.Range("I66").Value = 3376.6667
$ws.Range("J66").Value = 5086.75
$ws.Range("K66").Value = 16883.3335
$ws.Range("L66").Value = 25433.75
$ws.Range("M66").Value = -13451.3335
$ws.Range("N66").Value = -32297.75
$ws.Range("H67").Value = 28000
$ws.Range("J67").Value = 28000
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716
$ws.Range("H74").Value = 50001536
$ws.Range("I74").Value = 445.2
$ws.Range("J74").Value = 100002630
$ws.Range("K74").Value = 445.2
$ws.Range("L74").Value = 100002630
$ws.Range("M74").Value = 428.8
$ws.Range("N74").Value = -100004378
$ws.Range("H77").Value = 50001536
$ws.Range("I77").Value = 445.2
$ws.Range("J77").Value = 100002630
$ws.Range("K77").Value = 2226
$ws.Range("L77").Value = 500013150
$ws.Range("M77").Value = 2142
$ws.Range("N77").Value = -500021886

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2506.05
$ws.Range("I134").Value = 2424.7646
$ws.Range("J134").Value = 2966.6667
$ws.Range("K134").Value = 7274.293799999999
$ws.Range("L134").Value = 8900.000100000001
$ws.Range("M134").Value = -4739.293799999999
$ws.Range("N134").Value = -13970.0001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10587.975
$ws.Range("I31").Value = 1181.2
$ws.Range("K31").Value = 1181.2
$ws.Range("M31").Value = -886.2
$ws.Range("H34").Value = 10587.975
$ws.Range("I34").Value = 1181.2
$ws.Range("K34").Value = 1181.2
$ws.Range("M34").Value = -979.2
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 7709
$ws.Range("J49").Value = 7709
$ws.Range("L49").Value = 23127
$ws.Range("N49").Value = -23439
$ws.Range("H96").Value = 3675
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -19118
$ws.Range("H107").Value = 19608944
$ws.Range("I107").Value = 347.35294
$ws.Range("J107").Value = 29413242
$ws.Range("K107").Value = 1042.05882
$ws.Range("L107").Value = 88239726
$ws.Range("M107").Value = 877.94118
$ws.Range("N107").Value = -88243566
$ws.Range("H125").Value = 1847.3334
$ws.Range("I125").Value = 815
$ws.Range("J125").Value = 2142.2856
$ws.Range("K125").Value = 2445
$ws.Range("L125").Value = 6426.8568
$ws.Range("M125").Value = 2475
$ws.Range("N125").Value = -16266.8568
$ws.Range("H129").Value = 1685216.2
$ws.Range("I129").Value = 586.25
$ws.Range("J129").Value = 2166539
$ws.Range("K129").Value = 1758.75
$ws.Range("L129").Value = 6499617
$ws.Range("M129").Value = 3241.25
$ws.Range("N129").Value = -6509617
$ws.Range("H132").Value = 2385.9492
$ws.Range("I132").Value = 2330.8386
$ws.Range("J132").Value = 2446.9644
$ws.Range("K132").Value = 20977.5474
$ws.Range("L132").Value = 22022.6796
$ws.Range("M132").Value = -18447.5474
$ws.Range("N132").Value = -27082.6796

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 53600
$ws.Range("J140").Value = 53600
$ws.Range("L140").Value = 53600
$ws.Range("N140").Value = -63960

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1565.0769
$ws.Range("I16").Value = 1481
$ws.Range("J16").Value = 1637.1428
$ws.Range("K16").Value = 1481
$ws.Range("L16").Value = 1637.1428
$ws.Range("M16").Value = -1311
$ws.Range("N16").Value = -1977.1428
$ws.Range("H22").Value = 34334.332
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 51001.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 51001.5
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -51591.5
$ws.Range("H27").Value = 34334.332
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 51001.5
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 51001.5
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -51215.5
$ws.Range("H68").Value = 1908.2858
$ws.Range("I68").Value = 1762.5
$ws.Range("J68").Value = 1998
$ws.Range("K68").Value = 1762.5
$ws.Range("L68").Value = 1998
$ws.Range("M68").Value = -1013.5
$ws.Range("N68").Value = -3496
$ws.Range("H71").Value = 1908.2858
$ws.Range("I71").Value = 1762.5
$ws.Range("J71").Value = 1998
$ws.Range("K71").Value = 8812.5
$ws.Range("L71").Value = 9990
$ws.Range("M71").Value = -5068.5
$ws.Range("N71").Value = -17478

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 69864.5
$ws.Range("J138").Value = 69864.5
$ws.Range("L138").Value = 69864.5
$ws.Range("N138").Value = -80144.5
